$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    # Force the cell to be treated as literal text so that numeric-looking
    # or date-looking strings (e.g. "6481", "8/4/2025") are not silently
    # coerced into numbers / dates by Excel's type inference, then drop the
    # temporary "Text" number format so the cell keeps the default style.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# ---- Row 79 ----
Set-TextCell "A79" "6481"
Set-TextCell "B79" "8/4/2025"
$ws.Range("C79").Value = "DIAZ, CESAR, GRAL. 3520"
Set-TextCell "D79" "11"
Set-TextCell "E79" "808703875"
$ws.Range("F79").Value = "AYKO"
$ws.Range("G79").Value = "Pendiente"
$ws.Range("H79").Value = "Picada"
$ws.Range("I79").Value = 1
$ws.Range("J79").Value = "Cambio"
$ws.Range("K79").Value = "Sin equipos"
$ws.Range("L79").Value = "Pasante"
$ws.Range("M79").Value = -58.486002
$ws.Range("N79").Value = -34.61821
$ws.Range("O79").Value = "Devoto"
$ws.Range("P79").Value = "Capital Norte"

# ---- Row 80 ----
Set-TextCell "A80" "6556"
Set-TextCell "B80" "8/4/2025"
$ws.Range("C80").Value = "2 DE ABRIL DE 1982 6982"
Set-TextCell "D80" "8"
Set-TextCell "E80" ""
$ws.Range("F80").Value = "AYKO"
$ws.Range("G80").Value = "Pendiente"
$ws.Range("H80").Value = "Traspasar a columna o cortar redes en punta y desmontar poste"
$ws.Range("I80").Value = 1
$ws.Range("J80").Value = "Desmonte"
$ws.Range("K80").Value = "Sin equipos"
$ws.Range("L80").Value = "Poste"
$ws.Range("M80").Value = -58.494864
$ws.Range("N80").Value = -34.678826
$ws.Range("O80").Value = "Boedo"
$ws.Range("P80").Value = "Capital Sur"
